$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force text-number cells to remain stored as text (avoid numeric auto-conversion)
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D49").NumberFormat = "@"

# Apply updated values
$ws.Range("D2").Value = "42.628.07"
$ws.Range("E2").Value = "  -1.42%  "
$ws.Range("D3").Value = "2.366.90"
$ws.Range("E3").Value = "  +0.84%  "
$ws.Range("E4").Value = "  -0.19%  "
$ws.Range("D5").Value = "331.39"
$ws.Range("E5").Value = "  +5.82%  "
$ws.Range("D6").Value = "101.92"
$ws.Range("E6").Value = "  -6.79%  "
$ws.Range("E7").Value = "  -0.73%  "
$ws.Range("E8").Value = "  -0.01%  "
$ws.Range("D9").Value = "0.631"
$ws.Range("E9").Value = "  -1.15%  "
$ws.Range("D10").Value = "40.17"
$ws.Range("E10").Value = "  -6.57%  "
$ws.Range("D11").Value = "0.0922"
$ws.Range("E11").Value = "  -1.74%  "
$ws.Range("D12").Value = "8.49"
$ws.Range("E12").Value = "  -3.93%  "
$ws.Range("E13").Value = "  -3.05%  "
$ws.Range("E14").Value = "  +0.28%  "
$ws.Range("D15").Value = "16.55"
$ws.Range("E15").Value = "  +1.26%  "
$ws.Range("D16").Value = "2.726.15"
$ws.Range("E16").Value = "  +0.91%  "
$ws.Range("D17").Value = "2.363.99"
$ws.Range("E17").Value = "  +0.62%  "
$ws.Range("D18").Value = "8.20"
$ws.Range("E18").Value = "  +13.13%  "
$ws.Range("D19").Value = "42.738.50"
$ws.Range("E19").Value = "  -1.08%  "
$ws.Range("E20").Value = "  -1.67%  "
$ws.Range("D21").Value = "76.53"
$ws.Range("E21").Value = "  +1.50%  "
$ws.Range("E22").Value = "  +9.52%  "
$ws.Range("D23").Value = "270.82"
$ws.Range("E23").Value = "  +6.13%  "
$ws.Range("B24").Value = "InternetComputer(DFINITY)"
$ws.Range("C24").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D24").Value = "10.33"
$ws.Range("E24").Value = "  +13.63%  "
$ws.Range("B25").Value = "ImmutableX"
$ws.Range("C25").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D25").Value = "2.33"
$ws.Range("E25").Value = "  -9.70%  "
$ws.Range("E26").Value = "  -0.03%  "
$ws.Range("D27").Value = "11.56"
$ws.Range("E27").Value = "  -3.83%  "
$ws.Range("D28").Value = "23.26"
$ws.Range("E28").Value = "  +4.26%  "
$ws.Range("D29").Value = "2.21"
$ws.Range("E29").Value = "  -2.12%  "
$ws.Range("D30").Value = "176.80"
$ws.Range("E30").Value = "  +1.79%  "
$ws.Range("E31").Value = "  -2.31%  "
$ws.Range("D32").Value = "0.0904"
$ws.Range("E32").Value = "  -2.30%  "
$ws.Range("D33").Value = "35.60"
$ws.Range("E33").Value = "  -8.96%  "
$ws.Range("D34").Value = "6.13"
$ws.Range("E34").Value = "  +0.73%  "
$ws.Range("D35").Value = "0.133"
$ws.Range("E35").Value = "  +0.65%  "
$ws.Range("D36").Value = "4.64"
$ws.Range("E36").Value = "  -6.60%  "
$ws.Range("D37").Value = "2.99"
$ws.Range("E37").Value = "  +10.88%  "
$ws.Range("D38").Value = "0.0360"
$ws.Range("E38").Value = "  -4.54%  "
$ws.Range("E39").Value = "  +1.06%  "
$ws.Range("D40").Value = "3.82"
$ws.Range("E40").Value = "  -7.11%  "
$ws.Range("E41").Value = "  +3.39%  "
$ws.Range("D42").Value = "0.237"
$ws.Range("E42").Value = "  +1.67%  "
$ws.Range("D43").Value = "70.36"
$ws.Range("E43").Value = "  -3.15%  "
$ws.Range("E44").Value = "  -0.07%  "
$ws.Range("D45").Value = "119.96"
$ws.Range("E45").Value = "  +8.34%  "
$ws.Range("D46").Value = "92.23"
$ws.Range("E46").Value = "  +31.99%  "
$ws.Range("D47").Value = "11.88"
$ws.Range("E47").Value = "  -6.47%  "
$ws.Range("D49").Value = "9.22"
$ws.Range("E49").Value = "  -0.72%  "
$ws.Range("E50").Value = "  -2.71%  "
$ws.Range("D51").Value = "1.572.07"
$ws.Range("E51").Value = "  +5.24%  "

